# [24.06.06 17:30] Model auto-generation
# Rename sheet "# Type de vecteurs" -> "Type de vecteurs" and fix the
# "CISU" header value (and selection) on the "Type de ressource" and
# "Type de vecteurs" sheets.

$wb = $excel.ActiveWorkbook

# 1. Rename the "# Type de vecteurs" sheet (strip leading "# ").
$wsVecteurs = $wb.Worksheets.Item("# Type de vecteurs")
$wsVecteurs.Name = "Type de vecteurs"

$wsRessource = $wb.Worksheets.Item("Type de ressource")

# 2. "Type de vecteurs" sheet: fix B1 value + selection.
$wsVecteurs.Activate()
$wsVecteurs.Range("B1").Value = "CISU"
$wsVecteurs.Range("B4").Select()

# 3. "Type de ressource" sheet: fix B1 value + selection.
# Re-activate it last so it stays the active/selected tab, matching the
# original workbook state.
$wsRessource.Activate()
$wsRessource.Range("B1").Value = "CISU"
$wsRessource.Range("B2").Select()
